# Insert a new data row at row 42 (pushing existing rows 42..117 down to 43..118)
# and populate it with the new Ají price record for "Región del Maule".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A42").EntireRow.Insert()

$ws.Range("A42").Value = 7
$ws.Range("B42").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C42").Value = "Ñuble"
$ws.Range("D42").Value = 44914
$ws.Range("E42").Value = 16
$ws.Range("F42").Value = 100112021
$ws.Range("G42").Value = "Ají"
$ws.Range("H42").Value = "Americana (o)"
$ws.Range("I42").Value = "Primera"
$ws.Range("J42").Value = 80
$ws.Range("K42").Value = 13000
$ws.Range("L42").Value = 14000
$ws.Range("M42").Value = 13500
$ws.Range("N42").Value = "$/caja 15 kilos"
$ws.Range("O42").Value = "Región del Maule"
$ws.Range("P42").Value = 900
$ws.Range("Q42").Value = 15
$ws.Range("R42").Value = "Hortaliza"

# Match the date-cell number format used by other rows in column D
$ws.Range("D42").NumberFormat = $ws.Range("D43").NumberFormat
